$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add day15..day20 labels in S1:X1
$ws.Range("S1").Value = "day15"
$ws.Range("T1").Value = "day16"
$ws.Range("U1").Value = "day17"
$ws.Range("V1").Value = "day18"
$ws.Range("W1").Value = "day19"
$ws.Range("X1").Value = "day20"

# Rows 2-28: new survey data in columns Q-X (day15-day20)
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 0
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 1
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = 1
$ws.Range("U5").Value = 1
$ws.Range("V5").Value = 1
$ws.Range("W5").Value = 1
$ws.Range("X5").Value = 1
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 1
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 1
$ws.Range("U6").Value = 1
$ws.Range("V6").Value = 1
$ws.Range("W6").Value = 1
$ws.Range("X6").Value = 1
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 1
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 1
$ws.Range("U7").Value = 1
$ws.Range("V7").Value = 1
$ws.Range("W7").Value = 1
$ws.Range("X7").Value = 1
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 1
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 1
$ws.Range("U8").Value = 1
$ws.Range("V8").Value = 1
$ws.Range("W8").Value = 1
$ws.Range("X8").Value = 1
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 1
$ws.Range("S9").Value = 1
$ws.Range("T9").Value = 1
$ws.Range("U9").Value = 1
$ws.Range("V9").Value = 1
$ws.Range("W9").Value = 1
$ws.Range("X9").Value = 1
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 1
$ws.Range("S10").Value = 1
$ws.Range("T10").Value = 1
$ws.Range("U10").Value = 1
$ws.Range("V10").Value = 1
$ws.Range("W10").Value = 1
$ws.Range("X10").Value = 1
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 1
$ws.Range("S11").Value = 1
$ws.Range("T11").Value = 1
$ws.Range("U11").Value = 1
$ws.Range("V11").Value = 1
$ws.Range("W11").Value = 1
$ws.Range("X11").Value = 1
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 1
$ws.Range("S12").Value = 1
$ws.Range("T12").Value = 1
$ws.Range("U12").Value = 1
$ws.Range("V12").Value = 1
$ws.Range("W12").Value = 1
$ws.Range("X12").Value = 1
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 1
$ws.Range("S13").Value = 1
$ws.Range("T13").Value = 1
$ws.Range("U13").Value = 1
$ws.Range("V13").Value = 1
$ws.Range("W13").Value = 1
$ws.Range("X13").Value = 1
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = 1
$ws.Range("S14").Value = 1
$ws.Range("T14").Value = 1
$ws.Range("U14").Value = 1
$ws.Range("V14").Value = 1
$ws.Range("W14").Value = 1
$ws.Range("X14").Value = 1
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = 1
$ws.Range("S15").Value = 1
$ws.Range("T15").Value = 1
$ws.Range("U15").Value = 1
$ws.Range("V15").Value = 1
$ws.Range("W15").Value = 1
$ws.Range("X15").Value = 1
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = 1
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 1
$ws.Range("U16").Value = 1
$ws.Range("V16").Value = 1
$ws.Range("W16").Value = 0
$ws.Range("X16").Value = 0
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = 1
$ws.Range("S17").Value = 1
$ws.Range("T17").Value = 1
$ws.Range("U17").Value = 1
$ws.Range("V17").Value = 1
$ws.Range("W17").Value = 1
$ws.Range("X17").Value = 1
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = 1
$ws.Range("S18").Value = 1
$ws.Range("T18").Value = 1
$ws.Range("U18").Value = 1
$ws.Range("V18").Value = 1
$ws.Range("W18").Value = 1
$ws.Range("X18").Value = 1
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 1
$ws.Range("S19").Value = 1
$ws.Range("T19").Value = 1
$ws.Range("U19").Value = 1
$ws.Range("V19").Value = 1
$ws.Range("W19").Value = 0
$ws.Range("X19").Value = 0
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = 1
$ws.Range("S20").Value = 1
$ws.Range("T20").Value = 1
$ws.Range("U20").Value = 1
$ws.Range("V20").Value = 1
$ws.Range("W20").Value = 1
$ws.Range("X20").Value = 1
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = 1
$ws.Range("S21").Value = 1
$ws.Range("T21").Value = 0
$ws.Range("U21").Value = 0
$ws.Range("V21").Value = 0
$ws.Range("W21").Value = 0
$ws.Range("X21").Value = 0
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = 1
$ws.Range("S22").Value = 1
$ws.Range("T22").Value = 1
$ws.Range("U22").Value = 1
$ws.Range("V22").Value = 1
$ws.Range("W22").Value = 1
$ws.Range("X22").Value = 1
$ws.Range("Q23").Value = 0
$ws.Range("R23").Value = 0
$ws.Range("S23").Value = 0
$ws.Range("T23").Value = 0
$ws.Range("U23").Value = 0
$ws.Range("V23").Value = 0
$ws.Range("W23").Value = 0
$ws.Range("X23").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("R24").Value = 0
$ws.Range("S24").Value = 0
$ws.Range("T24").Value = 0
$ws.Range("U24").Value = 0
$ws.Range("V24").Value = 0
$ws.Range("W24").Value = 0
$ws.Range("X24").Value = 0
$ws.Range("Q25").Value = 0
$ws.Range("R25").Value = 0
$ws.Range("S25").Value = 0
$ws.Range("T25").Value = 0
$ws.Range("U25").Value = 0
$ws.Range("V25").Value = 0
$ws.Range("W25").Value = 0
$ws.Range("X25").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("R26").Value = 0
$ws.Range("S26").Value = 0
$ws.Range("T26").Value = 0
$ws.Range("U26").Value = 0
$ws.Range("V26").Value = 0
$ws.Range("W26").Value = 0
$ws.Range("X26").Value = 0
$ws.Range("Q27").Value = 0
$ws.Range("R27").Value = 0
$ws.Range("S27").Value = 0
$ws.Range("T27").Value = 0
$ws.Range("U27").Value = 0
$ws.Range("V27").Value = 0
$ws.Range("W27").Value = 0
$ws.Range("X27").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("R28").Value = 0
$ws.Range("S28").Value = 0
$ws.Range("T28").Value = 0
$ws.Range("U28").Value = 0
$ws.Range("V28").Value = 0
$ws.Range("W28").Value = 0
$ws.Range("X28").Value = 0

# Rows 29-36: new survey data in columns O-V (day15-day20, offset since fewer pre-existing day columns)
$ws.Range("O29").Value = 1
$ws.Range("P29").Value = 1
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = 1
$ws.Range("S29").Value = 1
$ws.Range("T29").Value = 1
$ws.Range("U29").Value = 1
$ws.Range("V29").Value = 1
$ws.Range("O30").Value = 1
$ws.Range("P30").Value = 1
$ws.Range("Q30").Value = 1
$ws.Range("R30").Value = 1
$ws.Range("S30").Value = 1
$ws.Range("T30").Value = 1
$ws.Range("U30").Value = 1
$ws.Range("V30").Value = 1
$ws.Range("O31").Value = 1
$ws.Range("P31").Value = 1
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = 1
$ws.Range("S31").Value = 1
$ws.Range("T31").Value = 1
$ws.Range("U31").Value = 1
$ws.Range("V31").Value = 1
$ws.Range("O32").Value = 1
$ws.Range("P32").Value = 1
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = 1
$ws.Range("S32").Value = 1
$ws.Range("T32").Value = 1
$ws.Range("U32").Value = 1
$ws.Range("V32").Value = 1
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = 0
$ws.Range("S33").Value = 0
$ws.Range("T33").Value = 0
$ws.Range("U33").Value = 0
$ws.Range("V33").Value = 0
$ws.Range("O34").Value = 1
$ws.Range("P34").Value = 1
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = 1
$ws.Range("S34").Value = 1
$ws.Range("T34").Value = 1
$ws.Range("U34").Value = 1
$ws.Range("V34").Value = 1
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = 0
$ws.Range("S35").Value = 0
$ws.Range("T35").Value = 0
$ws.Range("U35").Value = 0
$ws.Range("V35").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 0
$ws.Range("R36").Value = 0
$ws.Range("S36").Value = 0
$ws.Range("T36").Value = 0
$ws.Range("U36").Value = 0
$ws.Range("V36").Value = 0

# Leave the view with row 19 selected (whole row), matching the editor state at save time
$ws.Activate()
$ws.Rows.Item(19).Select()
